$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# Locate the table shape ("Table 6") on the slide.
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.HasTable) {
        $shp = $cand
    }
}

$tbl = $shp.Table

# EMU per point used by the PowerPoint object model.
$emuPerPt = 12700

# Resize the overall table/graphic-frame height (10902462 EMU wide stays the
# same, only the height shrinks from 3220533 -> 2946213 EMU).
$shp.Height = 2946213 / $emuPerPt

# Redistribute the two column widths (total stays 10902462 EMU).
$tbl.Columns.Item(1).Width = 4389122 / $emuPerPt
$tbl.Columns.Item(2).Width = 6513340 / $emuPerPt

# Fill in the previously-empty "PPT Link" cell (row 4, column 2) with the
# repository URL, styled like the other link cells in the table.
$cell = $tbl.Cell(4, 2)
$tr = $cell.Shape.TextFrame.TextRange
$tr.Text = "https://github.com/ASHWITHA2202/AI-mini-project/blob/main/ASHU_MINI_PROJECT_AI%5B1%5D%20(1).pptx"
$tr.Font.Name = "Times New Roman"
$tr.Font.Color.ObjectThemeColor = [int][Microsoft.Office.Core.MsoThemeColorIndex]::msoThemeColorAccent1
